$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Part 2")

# Row 2: Max P(wait) 0.9 -> 0.2 ; Number of Servers 2 -> 3
$ws.Cells.Item(2, 2).Value = 0.2
$ws.Cells.Item(2, 4).Value = 3

# Row 3: E(S) 0.02380952380952381 -> 0.06666666666666667
$ws.Cells.Item(3, 4).Value = 0.06666666666666667

# Row 4: Arrival Rate 15.0 -> 10.0 ; E(N) 0.02380952380952381 -> 1.9999999999999998
$ws.Cells.Item(4, 2).Value = 10.0
$ws.Cells.Item(4, 4).Value = 1.9999999999999998

# Row 5: Service Rate 21.0 -> 5.0 ; add (empty) C5/D5 cells
$ws.Cells.Item(5, 2).Value = 5.0
$ws.Cells.Item(5, 3).Font.Bold = $false
$ws.Cells.Item(5, 4).Font.Bold = $false
